$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 0.7129117291503684
$ws.Range("F3").Value = 1.032834171127953
$ws.Range("G3").Value = 78679.45245659979

$ws.Range("E4").Value = 0.01258079422247095
$ws.Range("F4").Value = 0.0008580503590280352
$ws.Range("G4").Value = 3954.239200401186

$ws.Range("E5").Value = 0.6589167194289471
$ws.Range("F5").Value = 1.030699407367167
$ws.Range("G5").Value = 66390.73424692501

$ws.Range("E6").Value = 0.7047538164791789
$ws.Range("F6").Value = 1.03222575368207
$ws.Range("G6").Value = 75895.06285653234

$ws.Range("E7").Value = 0.7128540410162563
$ws.Range("F7").Value = 1.032739146104849
$ws.Range("G7").Value = 78564.74524167427

$ws.Range("E8").Value = 0.7212907360074411
$ws.Range("F8").Value = 1.033354254656776
$ws.Range("G8").Value = 81262.6366782485

$ws.Range("E9").Value = 0.7662807552571944
$ws.Range("F9").Value = 1.037888771934041
$ws.Range("G9").Value = 94500.41535620989
